# "Generate Report for Archive"
#
# 1) The handoff status text moves from "Ready for handoff" to
#    "In Translation" everywhere it is used:
#      - Overview!E2 (zh-cn status) and Overview!F2 (de-de status)
#      - zh-cn!C2 (Status column)
#      - de-de!C2 (Status column)
#
# 2) The (now-narrower) status columns are resized down from
#    ~17.22 chars to ~13.41 chars on all three sheets:
#      - Overview columns E and F
#      - zh-cn column C
#      - de-de column C

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Update the status text ---
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value     = "In Translation"
$dede.Range("C2").Value     = "In Translation"

# --- Narrow the status columns ---
$overview.Range("E1:F1").EntireColumn.ColumnWidth = 12.5
$zhcn.Range("C1").EntireColumn.ColumnWidth = 12.5
$dede.Range("C1").EntireColumn.ColumnWidth = 12.5
